# StudentProfile.xlsx: refresh the sheet with a new student's data.
# Order of writes matters: Excel appends new shared-string entries in the
# order their owning cells are (re)written, so we touch A9 first, then
# B3, B8, B4 to reproduce the exact shared-string ordering of the target
# file (new entries 77-80: "Targetted Institution (e.g MIT)", "Chicão",
# "Senior", "Science Mater Tow").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A9: add the "(e.g MIT)" example hint (italic run) after the existing
# "Targetted Institution " label.
$cellA9 = $ws.Range("A9")
$prefix = "Targetted Institution "
$suffix = "(e.g MIT)"
$cellA9.Value = $prefix + $suffix
$cellA9.Characters($prefix.Length + 1, $suffix.Length).Font.Italic = $true

# B3: First Name
$ws.Range("B3").Value = "Chicão"

# B8: High School Year attended -> Senior
$ws.Range("B8").Value = "Senior"

# B4: Last Name
$ws.Range("B4").Value = "Science Mater Tow"

# B7: School Year attended -> 2018/2019
$ws.Range("B7").Value = "2018/2019"

# Leave the active selection on B4, matching the saved file.
$ws.Range("B4").Select()

$wb.Save()
